$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-run of the login test: rows 2-5 now report ERROR (row 6 remains SKIPPED)
$ws.Range("D2").Value = "ERROR"
$ws.Range("D3").Value = "ERROR"
$ws.Range("D4").Value = "ERROR"
$ws.Range("D5").Value = "ERROR"

# Apply the date-style number format used when the report was regenerated
$ws.Range("A1:D1").NumberFormat = "d-mmm-yy"
$ws.Range("A2:D6").NumberFormat = "d-mmm-yy"

# Scroll/selection state left behind after the last edit
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("A1:D6").Select()
$ws.Range("A2").Activate()
